# Update read number and md5 hash soda_rerun
# Fills in the previously-blank "Forward reads count", "Reverse reads
# count", "Forward MD5 checksum" and "Reverse MD5 checksum" columns
# (AB:AE) for samples in rows 3-10 of the SUMMARY.csv sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB3").Value = 67169028
$ws.Range("AC3").Value = 67169028
$ws.Range("AD3").Value = "a7b04c567045ef4ef8d3838cffb4d5a9"
$ws.Range("AE3").Value = "12c9fc0cdae4d667083c7f35fa7585b5"

$ws.Range("AB4").Value = 37585237
$ws.Range("AC4").Value = 37585237
$ws.Range("AD4").Value = "aab80aef28c4e1b0b971cb84af62d41b"
$ws.Range("AE4").Value = "5ab4dbfbf98d4f38febc320ad309f455"

$ws.Range("AB5").Value = 73593472
$ws.Range("AC5").Value = 73593472
$ws.Range("AD5").Value = "307e8b93e6141c5e1148331562142a68"
$ws.Range("AE5").Value = "df41ab47eee41cd52303fb82e4dfc446"

$ws.Range("AB6").Value = 45978734
$ws.Range("AC6").Value = 45978734
$ws.Range("AD6").Value = "dbd0c75cfd24bdb2271270c51be9b801"
$ws.Range("AE6").Value = "50b3eed03ca51dad6a874c577c1a1913"

$ws.Range("AB7").Value = 124471218
$ws.Range("AC7").Value = 124471218
$ws.Range("AD7").Value = "1cba9f85ccfa9bb083798b960143a79a"
$ws.Range("AE7").Value = "23cd46ecc5b984f862ad8867bb8d1b92"

$ws.Range("AB8").Value = 37770014
$ws.Range("AC8").Value = 37770014
$ws.Range("AD8").Value = "74c4f06936290d48d8ebcf99def294e6"
$ws.Range("AE8").Value = "a6b5acac07e4e09c4a64bd23641cf87a"

$ws.Range("AB9").Value = 58650577
$ws.Range("AC9").Value = 58650577
$ws.Range("AD9").Value = "55fd26f0c0e3964415ea8ae83a699454"
$ws.Range("AE9").Value = "626c5cb199a26a6c20a8b34fd3ada4ac"

$ws.Range("AB10").Value = 51786260
$ws.Range("AC10").Value = 51786260
$ws.Range("AD10").Value = "1d6470067b9aa22d7e0b2b70d9033737"
$ws.Range("AE10").Value = "1be86e35b31b69037545bbf7417e8811"

# Move the active cell / selection the way the author left it.
$ws.Range("AB15").Select()
